# Add a new "Haemophilus influenzae type b" row (row 15) to the "adults" sheet,
# plus a trailing blank-but-styled row (row 16), matching the upstream upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 15: new vaccine entry -------------------------------------------------

# A15 needs a brand-new font (Helvetica, theme color 1) that doesn't exist yet in
# this workbook's style table. Grab the Helvetica family from an existing styled
# cell (V13, cellXfs #12) via a format-only paste so we don't mint a throwaway
# font for "Helvetica + default color" before nudging just the color to theme 1.
$ws.Range("V13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Haemophilus influenzae type b "
$ws.Range("A15").Font.ThemeColor = 1

$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "Conditional"
$ws.Range("D15").Value = "See Conditions and Alternate Dosing"
$ws.Range("E15").Value = "X"
$ws.Range("F15").Value = "X"
$ws.Range("G15").Value = "X"
$ws.Range("H15").Value = "X"
$ws.Range("I15").Value = 6935
$ws.Range("J15").Value = 364635

# U15 reuses the existing "condition" style (cellXfs #8, seen on W14/Y14).
$ws.Range("W14").Copy()
$ws.Range("U15").PasteSpecial(-4122)
$ws.Range("U15").Value = "Anatomical or functional asplenia (including sickle cell disease), Hematopoietic stem cell transplant (HSCT)"

# W15 is plain/unstyled text.
$ws.Range("W15").Value = "Anatomical or functional asplenia (including sickle cell disease)"

# X15:AB15 reuse the existing "alternate dosing" style (cellXfs #12, seen on V13).
$ws.Range("V13").Copy()
$ws.Range("X15:AB15").PasteSpecial(-4122)
$ws.Range("X15").Value = "1 dose if previously did not receive Hib; if elective splenectomy, 1 dose preferably at least 14 days before splenectomy "
$ws.Range("Y15").Value = "Hematopoietic stem cell transplant (HSCT)"
$ws.Range("Z15").Value = "3-dose series 4 weeks apart starting 6–12 months after successful transplant, regardless of Hib vaccination history"
$ws.Range("AA15").Value = "Hematopoietic stem cell transplant (HSCT): "
$ws.Range("AB15").Value = "3-dose series 4 weeks apart starting 6–12 months after successful transplant, regardless of Hib vaccination history"

# The row used to carry a single empty placeholder cell (AF15); it's gone now
# that the row has real content.
$ws.Range("AF15").Clear()

# Row 15 needs to grow to fit the new wrapped text.
$ws.Rows.Item(15).RowHeight = 176

# --- Row 16: trailing blank styled row -----------------------------------------

$ws.Range("W14").Copy()
$ws.Range("U16:Y16").PasteSpecial(-4122)

# Leave the active selection on Z15, matching the author's last-saved view.
$null = $ws.Range("Z15").Select()
